# Applies the proofreading corrections described in the commit
# "Word de errores corregido" to Errores.docx.
#
# Each correction is applied with Range.Find.Execute using a
# MatchCase=$true, whole-text literal search (no wildcards) so the
# replacement only touches the intended sentence.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# Item 1, bullet 1: "...composition Finder-> Warranty..." gets an extra
# space inserted before "Finder" when Word's grammar checker re-splits
# the run (gramStart/gramEnd around "Finder->").
Replace-Text "We have changed the composition Finder-> Warranty" "We have changed the composition  Finder-> Warranty"

# Item 1, bullet 3: "It was wrong the direction." -> "The direction was wrong."
Replace-Text "It was wrong the direction." "The direction was wrong."

# Item 1, bullet 4: "Remove attributes ... es scafolding." rewritten in English.
Replace-Text "Remove attributes sender and recipient of Message because se entiende con la navegabilidad y es scafolding." "We removed the attributes sender and recipient of Message because just with the navigability is understandable, which means that is scaffoding."

# Item 1, bullet 5 (Report-Complaint): add "the" before "multiplicity" and
# translate the justification to English.
Replace-Text "We have changed multiplicity Report-Complaint 0..*->0..* to 0..*-> 1 because un mismo report no puede estar en distintos complaints." "We have changed the multiplicity Report-Complaint 0..*->0..* to 0..*-> 1 because the same report can not be in different complaints."

# Item 1, bullet 6 (Complaint-FixUpTask): translate justification to English.
Replace-Text "no tiene sentido que la misma complaint esté en varias fixUpTasks." "the same report can not be in different fixUpTasks."

# Item 1, bullet 7 (FixUpTask-Application): translate justification to English.
Replace-Text "no tiene sentido que una aplication esté en varias fixUpTasks." "the same report can not be in different fixUpTasks."

# Item 1, bullet 8 (Report-Referee): translate justification to English.
Replace-Text "no tiene sentido que un report pueda tener varios referees." "it makes no sense that a report may have several referees."

# Item 1, bullet 9 (User Account-Actor): translate justification to English.
Replace-Text "no tiene sentido que un actor tenga varias users accounts." "it makes no sense that an actor has several users accounts."

# Item 1, bullet 10 (FixUpTask-Category): translate justification to English.
Replace-Text "una category puede tener 0 o varios FixUpTasks." "a category can have zero or many FixUpTasks."

# Item 1, bullet 11 (FixUpTask-Finder): translate justification to English.
Replace-Text "no tiene sentido que un FixUpTask tenga un solo Finder." "it makes no sense that a FixUpTask may have just one Finder."

# Item 1, bullet 12: "We have changed the multiplicity for scaffolding." ->
# "...because of scaffolding."
Replace-Text "We have changed the multiplicity for scaffolding." "We have changed the multiplicity because of scaffolding."

# Item 1, bullet 13: translate the final note to English.
Replace-Text "Navegabilidad Message – Actor no estaba especificada." "Navegability Message – Actor was not specified."
